$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2/A3 hold text values (e.g. " 71") that look numeric; force text format
# so Excel does not silently convert them to numbers.
$ws.Range("A2:A3").NumberFormat = "@"

# Row 2 updates
$ws.Range("A2").Value = " 71"
$ws.Range("B2").Value = 108
$ws.Range("C2").Value = 1.07
$ws.Range("D2").Value = 248.4
$ws.Range("E2").Value = 1.23
$ws.Range("F2").Value = 161.2
$ws.Range("I2").Value = 0.17
$ws.Range("J2").Value = 0.19
$ws.Range("K2").Value = 0.15
$ws.Range("L2").Value = 98
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 18
$ws.Range("O2").Value = 20
$ws.Range("P2").Value = 16

# Row 3 updates
$ws.Range("A3").Value = " 29"
